$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) MARCELO's (003641655) Saldo changes from 178835.58 to 128835.58
$ws.Range("C2").Value = 128835.58

# 2) Remove the stale CRISTINA row (004853901 / 37.7) that currently sits
#    between JULIO (002401479, 37.84) and ANDRE (004384131, 37.47).
#    Do this before inserting new rows further up so this row index stays valid.
$ws.Rows(110).Delete()

# 3) Insert a new row for CRISTINA (004853901, 68734.99) right above THIAGO
#    (005064129), i.e. before current row 4. Force column A to text so the
#    leading zeros in the account number survive (matches the other rows).
$ws.Rows(4).Insert()
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "004853901"
$ws.Range("B4").Value = "CRISTINA"
$ws.Range("C4").Value = 68734.99

# 4) Insert a new row for JOSE (005146441, 20000) right above RODRIGO
#    (004392159). After the previous insert, RODRIGO moved from row 5 to
#    row 6, so insert before row 6.
$ws.Rows(6).Insert()
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "005146441"
$ws.Range("B6").Value = "JOSE"
$ws.Range("C6").Value = 20000
